$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("Dollie Desai", 9825082629, "U2,902 Happy Excellencia, Vesu", "Apr 2025")
    ,@("pratik agrwal", 9825291307, "Vesu B201 rajhans cromona", "Apr 2025")
    ,@("Leher Thakur", 9624004792, "C-201, Aakash Evergreen, Vesu", "Apr 2025")
    ,@("Aakash Damani", 7434032409, "Hotel Shreeji Vatika, Piplod, Dumas Road", "Apr 2025")
    ,@("Maheshbhai", 8320999203, "Vesu", "Apr 2025")
    ,@("yash", 8469894517, "Vesu", "Apr 2025")
    ,@("mamta", 8000906977, "Vesu JOLLY RECIDANCY", "Apr 2025")
    ,@("Suman Kothari", 9924599987, "Jolly Residency, Vesu", "Apr 2025")
    ,@("Daljit", 9898216896, "the grand plaza shop no 16", "Apr 2025")
    ,@("k.p", 9012406659, "Vesu", "May 2025")
    ,@("gourav bagaria", 6291673385, "A/B3 green victory", "May 2025")
    ,@("yashvi SINGHAN", 9408120230, "STAR GALAXY G-325", "May 2025")
    ,@("neha churiwal raj tilak apartment Q206", 9924624544, "Vesu", "May 2025")
    ,@("nikita", 6358107635, "happy glorius G-6-301", "May 2025")
    ,@("ankita damani", 8866040899, "Happy Elanza, C1201", "May 2025")
    ,@("sapna kothari", 9265764632, "B401 sapphire court", "May 2025")
    ,@("kirti gupta", 8347583195, "vesu someshwar m 80", "May 2025")
    ,@("diana", 9712902353, "Vesu JOLLY I1203", "May 2025")
    ,@("kunj bihari", 9651191359, "Vesu", "May 2025")
    ,@("shiva suwalka", 8875481372, "Vesu aagam hotel nova", "May 2025")
    ,@("pinki", 8511661008, "Vesu", "May 2025")
    ,@("shiva suwalka", 7081797878, "Vesu aagam hotel nova 306", "May 2025")
    ,@("anupama", 9016242079, "jolly H601", "May 2025")
    ,@("prieet", 9712755116, "Vesu", "May 2025")
    ,@("rahul baid", 9664951965, "Vesu", "May 2025")
    ,@("gavtabhai narnave", 7875055627, "hotel sunrise AGAM", "May 2025")
    ,@("babita", 8780331212, "rajhans sinfoniya B805", "May 2025")
    ,@("chirag shah", 9227907979, "Vesu", "May 2025")
    ,@("riya", 9460291537, "kassel braun appartment-8C", "May 2025")
    ,@("saloni", 9687422277, "subh anclev A1004", "Jun 2025")
    ,@("mitesh bhai", 9537890657, "dumas", "Jun 2025")
    ,@("janita saluja", 9825107160, "the majestic F 1303", "Jun 2025")
    ,@("zota", 8980878542, "Vesu", "Jun 2025")
    ,@("mitushree soni", 9327089115, "althan A1202 marvella residancy", "Jun 2025")
    ,@("renu", 9694966590, "happy glorios G6301", "Jun 2025")
    ,@("sakuntla", 9427344342, "hemtanpark b901", "Jun 2025")
    ,@("sneha juwani", 9638222737, "Fiona A401", "Jun 2025")
    ,@("vandna patel", 9106938198, "dream house102", "Jun 2025")
    ,@("salman salon", 9727182884, "Vesu cannon street 1", "Jun 2025")
    ,@("abhi", 8306540657, "d 301 sangini vedanta", "Jun 2025")
    ,@("jinali", 9555965750, "vesu", "Jun 2025")
    ,@("karishma murarka", 8160833677, "althan D 916 sentosa heights", "Jun 2025")
    ,@("kavita mehta", 9537547555, "The Atmoshphere, Vesu B201", "Jun 2025")
    ,@("Samir Bhai", 9971599922, "Vesu hemton park B301", "Jun 2025")
    ,@("harsh", 7047366468, "Vesu", "Jun 2025")
    ,@("vashar preet", 9898169655, "Vesu B2 501 celebrity greens", "Jun 2025")
    ,@("neha shah", 9724198940, "vesu", "Jun 2025")
    ,@("sourav", 9974173433, "jolly A1304", "Jun 2025")
    ,@("shradha jain", 8511031011, "dev raj residancy D801", "Jun 2025")
    ,@("amar singh", 9574802367, "Vesu", "Jun 2025")
    ,@("dr abhilasha", 7698173840, "Vesu", "Jun 2025")
    ,@("amar singh", 9512504076, "Vesu", "Jun 2025")
    ,@("deva park infinity", 7665340302, "ALTHAN", "Jul 2025")
    ,@("mihir", 8108349933, "celebrity greens D2 301", "Jul 2025")
    ,@("himansu", 9818580296, "swapna bhumi D603", "Jul 2025")
    ,@("jayshree vayeda", 8200812058, "Vesu shubhan B1004", "Jul 2025")
    ,@("aasha agrwal", 9327922291, "capital green", "Jul 2025")
    ,@("ishani", 9723423643, "hampton park B901", "Jul 2025")
    ,@("shuruchi", 9687618848, "401 maarquis hitsh", "Jul 2025")
    ,@("aagam", 9974842906, "celibrity green B2 501", "Jul 2025")
    ,@("foram", 7359724111, "Vesu", "Jul 2025")
    ,@("sureshbhai", 7080296220, "Vesu", "Jul 2025")
    ,@("falguni", 9586005322, "fiona B1102", "Jul 2025")
    ,@("mukesh", 9820019283, "Vesu", "Jul 2025")
    ,@("tanish jain", 7984925177, "green vally G415", "Jul 2025")
    ,@("Harsha", 9977937276, "Vesu", "Jul 2025")
    ,@("priyanshu", 9016717039, "maarqiuis hights 401", "Jul 2025")
    ,@("vanita", 9099921981, "Vesu", "Jul 2025")
    ,@("poja", 8780113488, "Vesu", "Jul 2025")
    ,@("nikunj", 9978062856, "Vesu", "Jul 2025")
    ,@("suman", 9327955366, "dream house 203", "Jul 2025")
    ,@("rajesh shah", 9825118388, "Vesu", "Jul 2025")
    ,@("avi raj", 9664784636, "Vesu", "Jul 2025")
    ,@("nirali shah", 9833869913, "piplodd", "Jul 2025")
)

$startRow = 275
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 4).Style = "Normal"
}